$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CST_REQ_03's description gains the "...and can view the histoty of previous
# transactions" clause that used to live in its own CST_REQ_05 row.
$ws.Range("C5").Value = "`nThe client should have direct access to his/her different accounts and  can view the histoty of previous transactions `n"
$ws.Rows(5).RowHeight = 15.75

# The now-redundant CST_REQ_05 row ("the client can view the histoty of
# previous transactions") is removed outright. Deleting it shifts every row
# below up by one and shrinks the A3:A7 merged cell down to A3:A6.
$ws.Rows(7).Delete()

# Row deletion shifted the CST_REQ_06 / CST_REQ_07 labels (and their
# descriptions) up along with everything else, so what is now row 7 reads
# "CST_REQ_06" and row 8 reads "CST_REQ_07". Relabel them back to the correct
# sequential ids (CST_REQ_05, CST_REQ_06) for the requirement table.
$ws.Range("B7").Value = "CST_REQ_05"
$ws.Range("B8").Value = "CST_REQ_06"

# Match the saved cursor position from the edit.
$ws.Range("C6").Select()
